$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 160, shifting existing rows 160-243 down to 161-244
$ws.Rows.Item(160).Insert()

# Fill the new row 160 with the required values
$ws.Cells.Item(160, 1).Value = 8
$ws.Cells.Item(160, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44572
$ws.Cells.Item(160, 4).NumberFormat = $ws.Cells.Item(159, 4).NumberFormat
$ws.Cells.Item(160, 5).Value = 4
$ws.Cells.Item(160, 6).Value = 100112032
$ws.Cells.Item(160, 7).Value = "Zapallo italiano"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 400
$ws.Cells.Item(160, 11).Value = 8000
$ws.Cells.Item(160, 12).Value = 9000
$ws.Cells.Item(160, 13).Value = 8500
$ws.Cells.Item(160, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(160, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(160, 16).Value = 121
$ws.Cells.Item(160, 17).Value = 70
$ws.Cells.Item(160, 18).Value = "Hortaliza"
